$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '75.763.54'
$ws.Range('E2').Value = '  +1.95%  '

# Row 3
$ws.Range('D3').Value = '2.841.07'
$ws.Range('E3').Value = '  +6.76%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '193.51'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.16%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '597.57'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.46%  '

# Row 7
$ws.Range('E7').Value = '  -0.08%  '

# Row 8
$ws.Range('E8').Value = '  +3.52%  '

# Row 9
$ws.Range('E9').Value = '  +0.73%  '

# Row 10
$ws.Range('D10').Value = '2.838.81'
$ws.Range('E10').Value = '  +6.75%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.388'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +9.52%  '

# Row 12
$ws.Range('E12').Value = '  -2.16%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.92'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +4.56%  '

# Row 14
$ws.Range('D14').Value = '3.357.25'
$ws.Range('E14').Value = '  +6.62%  '

# Row 15
$ws.Range('D15').Value = '75.643.08'
$ws.Range('E15').Value = '  +1.80%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.52'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.81%  '

# Row 17
$ws.Range('E17').Value = '  +1.94%  '

# Row 18
$ws.Range('D18').Value = '2.841.87'
$ws.Range('E18').Value = '  +6.60%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.09'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.48%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.44'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +4.64%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '383.64'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.73%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.31'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.52%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.14'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.05%  '

# Row 24
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.19'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.39%  '

# Row 25
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.05%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.25'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +4.07%  '

# Row 27
$ws.Range('D27').Value = '2.980.65'
$ws.Range('E27').Value = '  +6.50%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.75'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +5.16%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000104'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +11.80%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.23%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.43'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.13%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '518.27'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.24%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.74'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.10%  '

# Row 34
$ws.Range('E34').Value = '  +4.58%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.14%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '165.05'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.85%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.95'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +4.51%  '

# Row 38
$ws.Range('E38').Value = '  +0.91%  '

# Row 39
$ws.Range('E39').Value = '  +0.37%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '184.71'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +9.19%  '

# Row 41
$ws.Range('E41').Value = '  -0.03%  '

# Row 42
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.10'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.85%  '

# Row 43
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.344'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +5.91%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.69'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.48%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.23'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +4.26%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.09'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.71%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0883'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +4.92%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.38'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.75%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.572'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +8.89%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.76'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.90%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.655'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +10.97%  '
